$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 data
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1330808"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1330808"
$ws.Range("C2").Value = "PR & Editorial Assistant"
$ws.Range("D2").Value = "Londres, Royaume-Uni"
$ws.Range("F2").Value = "10 applicants"
$ws.Range("G2").Value = "6 - 18 Months"
$ws.Range("H2").Value = "Orenda books"

# Update row 3 data
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "1329595"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1329595"
$ws.Range("C3").Value = "Sales Intern"
$ws.Range("D3").Value = "Adana, Reşatbey, Seyhan/Adana, Türkiye"
$ws.Range("F3").Value = "63 applicants"
$ws.Range("G3").Value = "6 - 18 Months"
$ws.Range("H3").Value = "PROPER PERLİT"

# Delete rows 4-6 (the old data no longer present)
$ws.Range("A4:H6").EntireRow.Delete()

# Update column widths (ColumnWidth <-> stored width has a fixed offset in this engine)
$ws.Range("D1").ColumnWidth = 40.1666666666667
$ws.Range("H1").ColumnWidth = 15.1666666666667
